$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.268.14"
$ws.Range("E2").Value = "  +3.73%  "
$ws.Range("D3").Value = "1.607.73"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "'213.20"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "'0.485"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "'0.0618"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "'18.09"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").Value = "'0.0821"
$ws.Range("E11").Value = "  +5.17%  "
$ws.Range("D12").Value = "1.830.55"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").Value = "1.609.18"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "26.233.24"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").Value = "'60.71"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "'198.65"
$ws.Range("E20").Value = "  +7.17%  "
$ws.Range("D21").Value = "'4.26"
$ws.Range("E21").Value = "  +2.77%  "
$ws.Range("D22").Value = "'9.36"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").Value = "'142.68"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").Value = "'6.47"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E35").Value = "  +4.65%  "
$ws.Range("D36").Value = "1.106.79"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("D41").Value = "'0.500"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D43").Value = "1.742.84"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'92.85"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.11"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "0.0₆0113"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("E47").Value = "  +9.04%  "
$ws.Range("D48").Value = "'53.55"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'0.410"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("E51").Value = "  -0.42%  "
